$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-12-24 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-12-25 Wednesday", 2)

$tbl = $d.Tables.Item(1)

# Row 1 (1-indexed)
$tbl.Cell(1, 1).Range.Text = "240×8=1920"
$tbl.Cell(1, 2).Range.Text = "870×5=4350"
$tbl.Cell(1, 3).Range.Text = "971×7=6797"
$tbl.Cell(1, 4).Range.Text = "354×2=708"
$tbl.Cell(1, 5).Range.Text = "959×6=5754"

# Row 5
$tbl.Cell(5, 1).Range.Text = "425×9=3825"
$tbl.Cell(5, 2).Range.Text = "789×3=2367"
$tbl.Cell(5, 3).Range.Text = "233×2=466"
$tbl.Cell(5, 4).Range.Text = "913×7=6391"
$tbl.Cell(5, 5).Range.Text = "532×9=4788"

# Row 10
$tbl.Cell(10, 1).Range.Text = "740×3=2220"
$tbl.Cell(10, 2).Range.Text = "126×8=1008"
$tbl.Cell(10, 3).Range.Text = "201×6=1206"
$tbl.Cell(10, 4).Range.Text = "423×7=2961"
$tbl.Cell(10, 5).Range.Text = "639×7=4473"

# Row 15
$tbl.Cell(15, 1).Range.Text = "411×6=2466"
$tbl.Cell(15, 2).Range.Text = "886×6=5316"
$tbl.Cell(15, 3).Range.Text = "916×4=3664"
$tbl.Cell(15, 4).Range.Text = "811×4=3244"
$tbl.Cell(15, 5).Range.Text = "315×9=2835"

# Row 20
$tbl.Cell(20, 1).Range.Text = "690×9=6210"
$tbl.Cell(20, 2).Range.Text = "389×4=1556"
$tbl.Cell(20, 3).Range.Text = "571×4=2284"
$tbl.Cell(20, 4).Range.Text = "878×9=7902"
$tbl.Cell(20, 5).Range.Text = "513×3=1539"
